$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores every Price/Volume cell as text, even when it
# looks like a plain number (e.g. "247.21"). Force those particular cells
# to Text format first so the COM layer does not silently coerce the
# assigned string into a numeric value.
$textCells = @("D5", "D7", "D9", "D10", "D11", "D12", "D13", "D15", "D16", "D20", "D21", "D23", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.157.60'
$ws.Range("E2").Value = '  -1.67%  '
$ws.Range("D3").Value = '2.247.43'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '247.21'
$ws.Range("E5").Value = '  -2.02%  '
$ws.Range("E6").Value = '  -1.60%  '
$ws.Range("D7").Value = '74.48'
$ws.Range("E7").Value = '  -1.48%  '
$ws.Range("D9").Value = '0.617'
$ws.Range("E9").Value = '  -4.71%  '
$ws.Range("D10").Value = '41.27'
$ws.Range("E10").Value = '  +4.53%  '
$ws.Range("D11").Value = '0.0942'
$ws.Range("E11").Value = '  -4.30%  '
$ws.Range("D12").Value = '7.09'
$ws.Range("E12").Value = '  -5.79%  '
$ws.Range("D13").Value = '0.102'
$ws.Range("E13").Value = '  -3.76%  '
$ws.Range("D14").Value = '2.581.68'
$ws.Range("E14").Value = '  -1.92%  '
$ws.Range("D15").Value = '14.52'
$ws.Range("E15").Value = '  -4.20%  '
$ws.Range("D16").Value = '0.851'
$ws.Range("E16").Value = '  -2.47%  '
$ws.Range("D17").Value = '2.245.26'
$ws.Range("E17").Value = '  -1.45%  '
$ws.Range("D18").Value = '42.060.07'
$ws.Range("E18").Value = '  -1.66%  '
$ws.Range("D19").Value = '0.0₃0973'
$ws.Range("E19").Value = '  -2.99%  '
$ws.Range("D20").Value = '6.14'
$ws.Range("E20").Value = '  -1.70%  '
$ws.Range("D21").Value = '71.88'
$ws.Range("E21").Value = '  -0.66%  '
$ws.Range("E22").Value = '  +6.66%  '
$ws.Range("D23").Value = '229.99'
$ws.Range("E23").Value = '  -3.23%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = '11.06'
$ws.Range("E25").Value = '  -2.27%  '
$ws.Range("D26").Value = '3.55'
$ws.Range("E26").Value = '  -8.28%  '
$ws.Range("D27").Value = '7.72'
$ws.Range("E27").Value = '  +22.93%  '
$ws.Range("E28").Value = '  -4.37%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '169.60'
$ws.Range("E29").Value = '  +1.30%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '2.09'
$ws.Range("E30").Value = '  -3.70%  '
$ws.Range("D31").Value = '20.65'
$ws.Range("E31").Value = '  -1.95%  '
$ws.Range("D32").Value = '0.0827'
$ws.Range("E32").Value = '  -4.48%  '
$ws.Range("D33").Value = '0.119'
$ws.Range("E33").Value = '  -5.62%  '
$ws.Range("D34").Value = '30.13'
$ws.Range("E34").Value = '  -2.83%  '
$ws.Range("D35").Value = '0.125'
$ws.Range("E35").Value = '  -3.14%  '
$ws.Range("D36").Value = '4.52'
$ws.Range("E36").Value = '  -1.49%  '
$ws.Range("D37").Value = '4.90'
$ws.Range("E37").Value = '  +2.17%  '
$ws.Range("E38").Value = '  -1.41%  '
$ws.Range("D39").Value = '13.49'
$ws.Range("E39").Value = '  -0.86%  '
$ws.Range("D40").Value = '2.18'
$ws.Range("E40").Value = '  -5.50%  '
$ws.Range("D41").Value = '5.79'
$ws.Range("E41").Value = '  -2.26%  '
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").Value = '61.65'
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '108.95'
$ws.Range("E43").Value = '  +3.92%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '0.203'
$ws.Range("E44").Value = '  -3.70%  '
$ws.Range("D45").Value = '8.67'
$ws.Range("E45").Value = '  -4.78%  '
$ws.Range("E46").Value = '  -1.07%  '
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("E48").Value = '  -3.87%  '
$ws.Range("E49").Value = '  -1.41%  '
$ws.Range("D50").Value = '2.27'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").Value = '2.70'
$ws.Range("E51").Value = '  -0.99%  '
